$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells we touch so that
# numeric-looking strings (e.g. '259.54', '1.000') are preserved
# verbatim as text instead of being parsed into numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (Coin / Link / Price / Volume columns).
$ws.Range("D2").Value = '26.357.22'
$ws.Range("E2").Value = '  -3.14%  '
$ws.Range("D3").Value = '1.832.86'
$ws.Range("E3").Value = '  -2.65%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '259.54'
$ws.Range("E5").Value = '  -7.76%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.5198'
$ws.Range("E7").Value = '  -2.00%  '
$ws.Range("D8").Value = '0.3223'
$ws.Range("E8").Value = '  -8.89%  '
$ws.Range("D9").Value = '0.06735'
$ws.Range("E9").Value = '  -4.21%  '
$ws.Range("D10").Value = '18.65'
$ws.Range("E10").Value = '  -8.48%  '
$ws.Range("D11").Value = '0.7654'
$ws.Range("E11").Value = '  -7.01%  '
$ws.Range("D12").Value = '0.07675'
$ws.Range("E12").Value = '  -1.92%  '
$ws.Range("D13").Value = '1.827.52'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = '88.88'
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D15").Value = '5.024'
$ws.Range("E15").Value = '  -3.45%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("D17").Value = '14.05'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").Value = '1.0000'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '0.000007891'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").Value = '26.386.73'
$ws.Range("E20").Value = '  -3.12%  '
$ws.Range("D21").Value = '2.089.16'
$ws.Range("E21").Value = '  -1.75%  '
$ws.Range("D22").Value = '4.546'
$ws.Range("E22").Value = '  -4.88%  '
$ws.Range("D23").Value = '9.430'
$ws.Range("E23").Value = '  -7.16%  '
$ws.Range("D24").Value = '5.919'
$ws.Range("D25").Value = '2.280'
$ws.Range("E25").Value = '  -5.29%  '
$ws.Range("D26").Value = '145.14'
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("D27").Value = '1.639'
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("D28").Value = '16.94'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("D29").Value = '111.24'
$ws.Range("E29").Value = '  -3.21%  '
$ws.Range("D30").Value = '4.200'
$ws.Range("E30").Value = '  -5.07%  '
$ws.Range("D31").Value = '4.135'
$ws.Range("E31").Value = '  -5.97%  '
$ws.Range("D32").Value = '0.08742'
$ws.Range("D33").Value = '0.04841'
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("D34").Value = '1.127'
$ws.Range("E34").Value = '  -5.07%  '
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("E36").Value = '  -9.15%  '
$ws.Range("D37").Value = '3.096'
$ws.Range("E37").Value = '  -6.65%  '
$ws.Range("E38").Value = '  -5.94%  '
$ws.Range("E39").Value = '  -7.93%  '
$ws.Range("D40").Value = '0.4916'
$ws.Range("E40").Value = '  -7.44%  '
$ws.Range("D41").Value = '112.41'
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("D42").Value = '0.8891'
$ws.Range("E42").Value = '  -8.73%  '
$ws.Range("D43").Value = '6.123'
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("D44").Value = '0.9996'
$ws.Range("D45").Value = '7.707'
$ws.Range("E45").Value = '  -6.37%  '
$ws.Range("D46").Value = '0.4201'
$ws.Range("E46").Value = '  -9.00%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.106'
$ws.Range("E47").Value = '  -3.72%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1255'
$ws.Range("E48").Value = '  -8.58%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05872'
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").Value = '35.33'
$ws.Range("E50").Value = '  -3.84%  '
$ws.Range("D51").Value = '59.15'
$ws.Range("E51").Value = '  -4.26%  '
